# Applies the changes described by the diff:
#  - Sheet "Veicolo": add header "id" in A1 and value 253497 in A2, set column A width to 30
#  - Sheet "Release Date (RD)": replace row 2 data, delete rows 3 and 4
#  - Sheet "RD Tassative": replace row 2 data, delete rows 3 through 16

$wb = $excel.ActiveWorkbook

# --- Sheet "Veicolo" ---
$ws1 = $wb.Worksheets.Item("Veicolo")
$ws1.Range("A1").Value = "id"
$ws1.Range("A2").Value = 253497
$ws1.Columns.Item(1).ColumnWidth = 29.14

# --- Sheet "Release Date (RD)" ---
$ws2 = $wb.Worksheets.Item("Release Date (RD)")
$ws2.Range("A2").Value = 253668
$ws2.Range("B2").Value = 45910.58333333334
$ws2.Range("C2").Value = 0
$ws2.Rows("3:4").Delete()

# --- Sheet "RD Tassative" ---
$ws3 = $wb.Worksheets.Item("RD Tassative")
$ws3.Range("A2").Value = 253295
$ws3.Range("B2").Value = 45911.58333333334
$ws3.Rows("3:16").Delete()
